$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.198.88"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3
$ws.Range("D3").Value = "3.887.82"
$ws.Range("E3").Value = "  -0.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'527.97"
$ws.Range("E5").Value = "  +8.92%  "

# Row 6
$ws.Range("D6").Value = "'142.04"
$ws.Range("E6").Value = "  -2.93%  "

# Row 7
$ws.Range("E7").Value = "  -2.09%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
$ws.Range("D9").Value = "'0.715"
$ws.Range("E9").Value = "  -2.82%  "

# Row 10
$ws.Range("E10").Value = "  +0.96%  "

# Row 11
$ws.Range("D11").Value = "'0.0000329"
$ws.Range("E11").Value = "  -4.89%  "

# Row 12
$ws.Range("D12").Value = "'41.82"
$ws.Range("E12").Value = "  -3.18%  "

# Row 13
$ws.Range("D13").Value = "4.506.25"
$ws.Range("E13").Value = "  -0.60%  "

# Row 14
$ws.Range("E14").Value = "  -5.27%  "

# Row 15
$ws.Range("D15").Value = "4.027.43"
$ws.Range("E15").Value = "  +3.10%  "

# Row 16
$ws.Range("E16").Value = "  +7.46%  "

# Row 17
$ws.Range("E17").Value = "  -0.63%  "

# Row 18
$ws.Range("D18").Value = "'13.75"
$ws.Range("E18").Value = "  -3.37%  "

# Row 19
$ws.Range("D19").Value = "'19.61"
$ws.Range("E19").Value = "  -3.18%  "

# Row 20
$ws.Range("D20").Value = "69.143.12"

# Row 21
$ws.Range("D21").Value = "'423.11"
$ws.Range("E21").Value = "  -1.81%  "

# Row 22
$ws.Range("E22").Value = "  -5.53%  "

# Row 23
$ws.Range("D23").Value = "'14.09"
$ws.Range("E23").Value = "  -6.59%  "

# Row 24
$ws.Range("D24").Value = "'87.64"
$ws.Range("E24").Value = "  -1.57%  "

# Row 25
$ws.Range("E25").Value = "  +8.99%  "

# Row 26
$ws.Range("D26").Value = "'11.62"
$ws.Range("E26").Value = "  -0.47%  "

# Row 27
$ws.Range("D27").Value = "'10.50"
$ws.Range("E27").Value = "  -6.12%  "

# Row 28
$ws.Range("D28").Value = "'36.03"
$ws.Range("E28").Value = "  -4.51%  "

# Row 29
$ws.Range("D29").Value = "'694.76"
$ws.Range("E29").Value = "  -2.86%  "

# Row 30
$ws.Range("D30").Value = "'13.07"
$ws.Range("E30").Value = "  -5.24%  "

# Row 32
$ws.Range("E32").Value = "  -4.10%  "

# Row 33
$ws.Range("D33").Value = "'67.84"
$ws.Range("E33").Value = "  +11.13%  "

# Row 34
$ws.Range("D34").Value = "'0.440"
$ws.Range("E34").Value = "  +9.32%  "

# Row 35
$ws.Range("D35").Value = "'5.89"
$ws.Range("E35").Value = "  -5.39%  "

# Row 36
$ws.Range("D36").Value = "'39.91"
$ws.Range("E36").Value = "  -4.25%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0830"
$ws.Range("E37").Value = "  -7.46%  "

# Row 38
$ws.Range("E38").Value = "  +3.07%  "

# Row 39
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.16%  "

# Row 40
$ws.Range("E40").Value = "  -0.15%  "

# Row 41
$ws.Range("D41").Value = "'0.0478"
$ws.Range("E41").Value = "  -2.18%  "

# Row 42
$ws.Range("D42").Value = "'2.76"
$ws.Range("E42").Value = "  -9.33%  "

# Row 43
$ws.Range("D43").Value = "'3.00"
$ws.Range("E43").Value = "  +0.46%  "

# Row 44
$ws.Range("D44").Value = "'2.94"
$ws.Range("E44").Value = "  -6.11%  "

# Row 45
$ws.Range("D45").Value = "'3.33"
$ws.Range("E45").Value = "  -0.84%  "

# Row 46
$ws.Range("E46").Value = "  -2.19%  "

# Row 47
$ws.Range("D47").Value = "'3.02"
$ws.Range("E47").Value = "  +7.72%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'26.47"
$ws.Range("E48").Value = "  +4.65%  "

# Row 49
$ws.Range("D49").Value = "'3.26"
$ws.Range("E49").Value = "  -4.90%  "

# Row 50
$ws.Range("D50").Value = "'142.42"
$ws.Range("E50").Value = "  -2.08%  "

# Row 51
$ws.Range("E51").Value = "  -4.23%  "

# Reset style for quote-prefixed numeric-text cells to avoid unintended style/format changes
foreach ($addr in @("D5","D6","D9","D11","D12","D18","D19","D21","D23","D24","D26","D27","D28","D29","D30","D33","D34","D35","D36","D39","D41","D42","D43","D44","D45","D47","D48","D49","D50")) {
    $ws.Range($addr).Style = "Normal"
}
